$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: Insert two new columns at D:E, shifting existing D:M data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: Copy cell formatting (number formats/fonts) from column F into new D:E columns
# so the new columns inherit the correct per-row style (date vs number format).
$ws.Range("F5:F102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Populate the new D and E columns with the new quarter figures.
$rowData = @"
7|43465|43373
8|2592200|2214700
9|2053500|1741000
10|538700|473700
11|BLANK|BLANK
12|88000|83300
13|0|0
14|1900|1500
15|15500|15300
16|BLANK|BLANK
17|2433100|2103400
18|159100|111300
19|BLANK|BLANK
20|-32400|-26100
21|197300|155500
22|0|0
23|126700|85200
24|45600|23900
25|0|0
26|81100|61300
27|90200|71100
28|0|0
29|8500|NA
30|0|0
31|0|0
32|32400|26100
33|98700|71100
34|0|0
35|98700|71100
38|43465|43373
39|BLANK|BLANK
40|BLANK|BLANK
41|326100|292700
42|0|0
43|880300|1003900
44|1908700|2101800
45|422300|390600
46|3537400|3789000
47|400000|419200
48|1373100|1367800
49|2068600|2085200
50|0|0
51|0|0
52|247300|255900
53|0|0
54|7626400|7917100
55|BLANK|BLANK
56|BLANK|BLANK
57|865900|855300
58|184200|186800
59|1716600|1617000
60|2766700|2659100
61|1275300|1699300
62|590900|589800
63|0|0
64|0|0
65|0|0
66|4693500|5011000
67|BLANK|BLANK
68|0|0
69|0|0
70|0|0
71|0|0
72|4477300|4405400
73|0|0
74|0|0
75|0|0
76|2932900|2906100
77|0|0
80|43465|43373
81|98700|71100
82|BLANK|BLANK
83|70600|70300
84|0|0
85|0|0
86|0|0
87|0|0
88|0|0
89|599900|200300
90|BLANK|BLANK
91|-64800|-48700
92|0|0
93|0|0
94|-64200|-48400
95|BLANK|BLANK
96|-11500|-11800
97|0|0
98|0|0
99|0|0
100|-505300|-128500
101|3000|-11300
102|33400|12100
"@

$lines = $rowData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|'
    $r = [int]$parts[0]
    $dVal = $parts[1]
    $eVal = $parts[2]

    if ($dVal -eq "BLANK") {
        # leave blank (already blank from copy/insert)
    } elseif ($dVal -eq "NA") {
        $ws.Cells.Item($r, 4).Value = "NA"
    } else {
        $ws.Cells.Item($r, 4).Value = [double]$dVal
    }

    if ($eVal -eq "BLANK") {
        # leave blank
    } elseif ($eVal -eq "NA") {
        $ws.Cells.Item($r, 5).Value = "NA"
    } else {
        $ws.Cells.Item($r, 5).Value = [double]$eVal
    }
}
